$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 6. This shifts the existing data rows
# 6-17 down to 7-18 (preserving their values and formatting, e.g. the date
# number format on column D), matching the diff which re-numbers every
# existing record down by one row and extends the used range to A1:T18.
$ws.Rows(6).Insert()

# Populate the newly inserted row 6 with the new record's data.
$ws.Range("A6").Value = 10
$ws.Range("B6").Value = "Vega Modelo de Temuco"
$ws.Range("C6").Value = "La Araucanía"
$ws.Range("D6").Value = 44901
$ws.Range("E6").Value = 9
$ws.Range("F6").Value = "Fruta"
$ws.Range("G6").Value = 100104
$ws.Range("H6").Value = "Frutos de pepita"
$ws.Range("I6").Value = 100104004
$ws.Range("J6").Value = "Níspero"
$ws.Range("K6").Value = "Californiana(o)"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 40
$ws.Range("N6").Value = 25000
$ws.Range("O6").Value = 25000
$ws.Range("P6").Value = 25000
$ws.Range("Q6").Value = "$/bandeja 10 kilos"
$ws.Range("R6").Value = "Provincia de Quillota"
$ws.Range("S6").Value = 2500
$ws.Range("T6").Value = 10
